$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: collapse the three runs ("Sprint 2 " / "startup /" /
#    " standup") -- split apart only because of the gramStart/gramEnd
#    proofing-error markers -- into a single run, dropping the proofErr
#    markers in the process.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleTextRange = $d.Range($titleRange.Start, $titleRange.End - 1)
# Force a real text mutation (setting identical text is treated as a no-op
# by the engine), then set the final wording.
$titleTextRange.Text = "zzz_placeholder_zzz"
$titlePara2 = $d.Paragraphs(1)
$titleRange2 = $titlePara2.Range
$titleTextRange2 = $d.Range($titleRange2.Start, $titleRange2.End - 1)
$titleTextRange2.Text = "Sprint 2 startup / standup"

# ---------------------------------------------------------------------------
# 2) Fill in the empty "Jack: " bullet under "Afgelopen tijd gedaan:" (06/03
#    standup) with the note about the "over bedrijf" landing page.
# ---------------------------------------------------------------------------
$jackDoneBullet = $d.Paragraphs(82)
$jackDoneBullet.Range.InsertAfter("Tussen pagina “over bedrijf” gemaakt. Zo hoeven wij niet 3 extra knoppen in de navbar toe te voegen.")
$jackDoneBullet.Range.Font.Size = 12
$jackDoneBullet.Range.Font.SizeBi = 12

# ---------------------------------------------------------------------------
# 3) "Komende tijd te doen" > Jack: replace the lone placeholder "A" bullet
#    with real text, then append a second bullet after it.
# ---------------------------------------------------------------------------
$jackTodoBullet = $d.Paragraphs(91)
$jackTodoRange = $jackTodoBullet.Range
$jackTodoTextRange = $d.Range($jackTodoRange.Start, $jackTodoRange.End - 1)
$jackTodoTextRange.Text = "Ecovriendelijkheid pagina afronden"

$jackTodoBullet2 = $d.Paragraphs(91)
$jackTodoBullet2.Range.InsertParagraphAfter()
$newTodoBullet = $d.Paragraphs(92)
$newTodoBullet.Range.InsertAfter("Product klachten formulier maken")
$newTodoBullet.Range.Font.Size = 12
$newTodoBullet.Range.Font.SizeBi = 12

# ---------------------------------------------------------------------------
# 4) "Potentiële obstakels" > Jack: replace the lone placeholder "A" bullet
#    with real text, then append a second bullet (with a
#    lastRenderedPageBreak marker, as in the source) after it.
#    (Index +1 vs. the original document because step 3 above already
#    inserted one extra paragraph earlier in the story.)
# ---------------------------------------------------------------------------
$jackObstacleBullet = $d.Paragraphs(99)
$jackObstacleRange = $jackObstacleBullet.Range
$jackObstacleTextRange = $d.Range($jackObstacleRange.Start, $jackObstacleRange.End - 1)
$jackObstacleTextRange.Text = "Bij ecovriendelijkheid pagina nog niet helemaal design bedacht, dus daar mee verden gaan."

$jackObstacleBullet2 = $d.Paragraphs(99)
$jackObstacleBullet2.Range.InsertParagraphAfter()
$newObstacleBullet = $d.Paragraphs(100)
$newObstacleXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val='24'/><w:szCs w:val='24'/></w:rPr><w:lastRenderedPageBreak/><w:t>Voor het eerst in het project werken met php, databases en formulieren, ik weet hoe  het werkt maar toch even extra aandacht aan moeten besteden.</w:t></w:r></w:p>"
$newObstacleBullet.Range.InsertXML($newObstacleXml)

Write-Output "done"
